$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (and implicitly the workbook's sheet tab) to reflect new "through" date
$ws.Name = "Through 2021-09-23"

# Row 10 (August) - update 2021 columns (T,U,V)
$ws.Range("U10").Value = 154
$ws.Range("V10").Value = 0.0314

# Row 11 (September) - label + most columns updated with new arrest data
$ws.Range("A11").Value = "September (through 09-23)"
$ws.Range("C11").Value = 24
$ws.Range("D11").Value = 0.04
$ws.Range("I11").Value = 53
$ws.Range("J11").Value = 0.0702
$ws.Range("L11").Value = 42
$ws.Range("M11").Value = 0.08699999999999999
$ws.Range("N11").Value = 5
$ws.Range("O11").Value = 52
$ws.Range("P11").Value = 0.0877
$ws.Range("R11").Value = 84
$ws.Range("S11").Value = 0.0345
$ws.Range("U11").Value = 135

# Row 12 (Total) - recompute totals impacted by the above changes
$ws.Range("C12").Value = 189
$ws.Range("D12").Value = 0.137
$ws.Range("I12").Value = 559
$ws.Range("J12").Value = 0.0806
$ws.Range("L12").Value = 475
$ws.Range("M12").Value = 0.1138
$ws.Range("N12").Value = 41
$ws.Range("O12").Value = 365
$ws.Range("P12").Value = 0.101
$ws.Range("R12").Value = 820
$ws.Range("S12").Value = 0.0586
$ws.Range("U12").Value = 1133
$ws.Range("V12").Value = 0.0605
